$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# columns that change, in order: B, D, E, F, G, J, K, O -> col indices 2,4,5,6,7,10,11,15
$colIndices = @(2, 4, 5, 6, 7, 10, 11, 15)

$values = @(
    @(2, 0.1423567229835641, 0.1938316955690311, 0.1678188689853357, 1.424625665348891, 0.002461622095605746, 0.1876083089544167, 0.947100573290129, 3.433484545394862),
    @(3, 0.1329742956853437, 0.1875410187144126, 0.1635095156113664, 1.429586174842939, 0.00246457404570447, 0.1836838837240222, 0.8357129766774847, 3.460978184578011),
    @(4, 0.1272837297545806, 0.1837423952983528, 0.1609436858420281, 1.433551216286176, 0.002466483008553415, 0.1813847363257821, 0.767163545317203, 3.480227812198393),
    @(5, 0.1249825994363505, 0.1822105931629778, 0.1599182770510765, 1.435397988586118, 0.002467285253904028, 0.1804755842034282, 0.7391912291239464, 3.488667023374163),
    @(6, 0.1246015792504949, 0.1819572181879607, 0.1597492291221592, 1.435718589667133, 0.002467419937667348, 0.1803262971388335, 0.7345442080682574, 3.490104253020263),
    @(7, 0.1272526235500493, 0.1837216712751513, 0.1609297750300591, 1.433575187482106, 0.002466493729263811, 0.1813723627567825, 0.7667864520496721, 3.480339218876395),
    @(8, 0.1391071606987566, 0.1916494730678124, 0.1663163939488399, 1.426145148612513, 0.002462619953329199, 0.186232233225553, 0.9087276431490636, 3.442472346218779),
    @(9, 0.1629063198143683, 0.2076990293270597, 0.177514395714816, 1.418878087841392, 0.002455785450148952, 0.1966400555362213, 1.185771313216264, 3.387042461601197),
    @(10, 0.1807232760364172, 0.2197938253163869, 0.1861284158575387, 1.418006162884993, 0.002451223942101092, 0.2048243792912814, 1.388465304218528, 3.357845946672541),
    @(11, 0.1888997423138363, 0.2253612365348943, 0.1901311913580628, 1.418582869061822, 0.002449247632411276, 0.2086650378802233, 1.480481119016076, 3.347077893924535),
    @(12, 0.1920061071342758, 0.2274787939983014, 0.1916590274269012, 1.418941454968362, 0.002448513378968945, 0.2101363302632535, 1.515296394022926, 3.343362630852596),
    @(13, 0.1913366485767085, 0.2270223282308734, 0.1913294442989155, 1.418857987944463, 0.002448670886194238, 0.2098187084844341, 1.507799623420112, 3.34414664841978),
    @(14, 0.189155103018436, 0.2255352635038719, 0.1902566456305479, 1.418609558993339, 0.002449186942177888, 0.2087857428522426, 1.483345985841083, 3.346764969436293),
    @(15, 0.1878201558604928, 0.2246256012914216, 0.1896010964541972, 1.41847565415317, 0.002449504879782276, 0.2081552252008265, 1.468363584827387, 3.348415984350652),
    @(16, 0.1801903487695711, 0.2194312891635946, 0.1858685164175355, 1.417988080018404, 0.002451355079611633, 0.204575750036696, 1.382447871713111, 3.358600322221804),
    @(17, 0.1755278945003482, 0.2162614269765584, 0.1836002416782634, 1.417938433426968, 0.002452515359953571, 0.2024099787338685, 1.329691310794544, 3.365492524654599),
    @(18, 0.1728529082276538, 0.2144443709338475, 0.1823035189824225, 1.418001480338674, 0.002453192020761772, 0.2011753497630764, 1.299329283116833, 3.369693273468016),
    @(19, 0.1719483652727263, 0.213830209121511, 0.1818658339238723, 1.418038552609488, 0.002453422725209251, 0.2007592257311757, 1.289046208540015, 3.371156172313562),
    @(20, 0.1760235249151805, 0.2165982269695661, 0.1838408830579041, 1.417934236199386, 0.002452390884207832, 0.202639383471805, 1.335309196424305, 3.364734350688195),
    @(21, 0.1897956025931222, 0.225971798935177, 0.1905714252806376, 1.418678721759591, 0.00244903498097201, 0.2090886906434122, 1.490529415614674, 3.34598606333185),
    @(22, 0.1988553327437188, 0.2321521255675094, 0.195040570159513, 1.419982585644306, 0.002446924045175743, 0.2134023042671487, 1.591804310767429, 3.335845536219949),
    @(23, 0.194014649141991, 0.2288486494915105, 0.1926488792688374, 1.419211829056138, 0.002448043179354284, 0.2110910198953206, 1.537768148803082, 3.34106412522496),
    @(24, 0.1757994331965023, 0.2164459430527472, 0.1837320662292896, 1.417935848484376, 0.002452447129849285, 0.2025356368797873, 1.332769448569138, 3.365076378810159),
    @(25, 0.1564093983431007, 0.2033037019243693, 0.1744170926792847, 1.420060613378041, 0.002457553280047555, 0.19373023148448, 1.110968955487976, 3.400016945167181)
)

foreach ($rowData in $values) {
    $r = $rowData[0]
    for ($i = 0; $i -lt $colIndices.Count; $i++) {
        $c = $colIndices[$i]
        $v = $rowData[$i + 1]
        $ws.Cells.Item($r, $c).Value = $v
    }
}

Write-Output "Updated $($values.Count) rows across columns B,D,E,F,G,J,K,O"
